# ch_hist.xlsx -- "some improvements to image floating / fixed some overfull
# boxes / updated sizes in histogram spreadsheets"
#
# 1. Shrink & reposition the embedded histogram chart (it previously spanned
#    out to col M / row 22, now it ends much sooner at col J / row 14 -- a
#    straightforward resize of the chart object).
# 2. Nudge the saved window y-position.
# 3. Leave the last active-cell selection on the sheet at J2 instead of O13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Resize the chart -----------------------------------------------
# The chart is anchored (two-cell anchor) from col D/row 3 and used to
# stretch to col M/row 22; shrink it so the bottom-right anchor lands on
# col J (index 9) / row 14 (index 13), matching the smaller histogram
# layout used elsewhere in the document.
$co = $ws.ChartObjects(1)
$co.Width = 346.0896
$co.Height = 170.07874

# --- 2. Window position --------------------------------------------------
# Saved workbook window moved down slightly on screen.
$excel.ActiveWindow.Top = 2400

# --- 3. Active selection ---------------------------------------------------
[void]$ws.Range("J2").Select()
